$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sheet is protected; unprotect before editing, then restore protection after
$ws.Unprotect()

# Update the confidentiality / as-of date notice text
$ws.Range("A9").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-03-19 for illustrative purposes only and are subject to change."

# Update Weight (D) and Percent Change (E) values for rows 2-6
$ws.Range("D2").Value = 0.250660239759725
$ws.Range("E2").Value = -0.004975124378109652

$ws.Range("D3").Value = 0.2524061765412907
$ws.Range("E3").Value = -0.01155401502021947

$ws.Range("D4").Value = 0.2464546212445252
$ws.Range("E4").Value = -0.003000230786983638

$ws.Range("D5").Value = 0.2504789624544592
$ws.Range("E5").Value = 0.008694489031567443

$ws.Range("E6").Value = -0.002725004774959183

# Re-protect the sheet to restore original protection state
$ws.Protect($null, $true, $true, $true)
